$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 69.08
$ws.Range("D2").Value = 7.69
$ws.Range("F2").Value = 12.62
$ws.Range("K2").Value = 0.77
$ws.Range("N2").Value = 0.92
$ws.Range("Q2").Value = 0.15
$ws.Range("U2").Value = 0.15
$ws.Range("V2").Value = 1.23
$ws.Range("C3").Value = 76
$ws.Range("D3").Value = 26
$ws.Range("F3").Value = 30
$ws.Range("C4").Value = 75.91
$ws.Range("D4").Value = 33.76
$ws.Range("F4").Value = 61.31
$ws.Range("K4").Value = 0.73
$ws.Range("M4").Value = 0.18
$ws.Range("N4").Value = 0.36
$ws.Range("P4").Value = 0.18
$ws.Range("Q4").Value = 1.82
$ws.Range("V4").Value = 3.1
$ws.Range("C5").Value = 82.75
$ws.Range("D5").Value = 66.75
$ws.Range("F5").Value = 32.75
$ws.Range("G5").Value = 0.75
$ws.Range("K5").Value = 0.5
$ws.Range("M5").Value = 0.5
$ws.Range("N5").Value = 0.25
$ws.Range("S5").Value = 0.25
$ws.Range("V5").Value = 2.75
$ws.Range("C6").Value = 57.76
$ws.Range("D6").Value = 56.71
$ws.Range("F6").Value = 14.08
$ws.Range("G6").Value = 0.26
$ws.Range("K6").Value = 0.13
$ws.Range("N6").Value = 0.39
$ws.Range("P6").Value = 0.26
$ws.Range("Q6").Value = 0.53
$ws.Range("U6").Value = 0.13
$ws.Range("V6").Value = 0.66
$ws.Range("C7").Value = 47.64
$ws.Range("D7").Value = 18.49
$ws.Range("F7").Value = 21.81
$ws.Range("K7").Value = 0.51
$ws.Range("N7").Value = 1.79
$ws.Range("P7").Value = 0.32
$ws.Range("Q7").Value = 0.51
$ws.Range("C8").Value = 52.96
$ws.Range("D8").Value = 31.48
$ws.Range("F8").Value = 44.07
$ws.Range("J8").Value = 5.19
$ws.Range("K8").Value = 0.37
$ws.Range("N8").Value = 0.37
$ws.Range("P8").Value = 0.74
$ws.Range("Q8").Value = 0.74
$ws.Range("V8").Value = 1.48
$ws.Range("C9").Value = 57.27
$ws.Range("D9").Value = 12.1
$ws.Range("F9").Value = 35.41
$ws.Range("G9").Value = 0.07000000000000001
$ws.Range("I9").Value = 0.07000000000000001
$ws.Range("P9").Value = 0.21
$ws.Range("R9").Value = 0.04
$ws.Range("U9").Value = 0.07000000000000001
$ws.Range("V9").Value = 0.28
$ws.Range("C10").Value = 61.56
$ws.Range("D10").Value = 14.21
$ws.Range("F10").Value = 29.21
$ws.Range("G10").Value = 0.16
$ws.Range("H10").Value = 2.25
$ws.Range("I10").Value = 0.08
$ws.Range("K10").Value = 0.48
$ws.Range("L10").Value = 0.08
$ws.Range("N10").Value = 0.72
$ws.Range("P10").Value = 0.32
$ws.Range("V10").Value = 0.5600000000000001
$ws.Range("C11").Value = 41.48
$ws.Range("D11").Value = 58.89
$ws.Range("F11").Value = 23.33
$ws.Range("H11").Value = 4.07
$ws.Range("J11").Value = 5.93
$ws.Range("K11").Value = 0.37
$ws.Range("N11").Value = 0.37
$ws.Range("Q11").Value = 1.48
$ws.Range("T11").Value = 15.93
$ws.Range("V11").Value = 1.11
